# 0.4.2 rare enemy quick-firing gun.
# Inserts a new mob row ("QuickFiringGun") above the Warlock row (row 25),
# pushing the existing rows 25-36 down to 26-37.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 25 - shifts old rows 25-36 down to 26-37.
$ws.Rows.Item(25).Insert()

# Populate the new row with the QuickFiringGun stats.
$ws.Range("A25").Value = "QuickFiringGun"
$ws.Range("B25").Value = 65
$ws.Range("C25").Value = 30
$ws.Range("D25").Value = 18
$ws.Range("E25").Value = 15
$ws.Range("F25").Value = 22
$ws.Range("G25").Value = 0.3
$ws.Range("I25").Value = 6
$ws.Range("J25").Value = 12
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 10
$ws.Range("O25").Value = "LIGHT:-0.2"
$ws.Range("P25").Value = "MACHINE"

# Match the author's last selection in the saved file.
[void]$ws.Range("B25").Select()
